# This document was edited so that consecutive runs carrying the same
# (or reconcilable) text were consolidated into single runs. The visible
# text does not change - only run/formatting boundaries collapse (e.g. the
# italicized species names lose their separate run and take on the
# surrounding run's formatting). We reproduce that by re-finding each
# affected span and replacing it with itself: Word's Find/Replace merges
# the matched range into a single run using the first run's formatting,
# exactly mirroring the target XML.

$d = $word.ActiveDocument

function Merge-Span($searchText) {
    $find = $d.Content.Find
    $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2) | Out-Null
}

# Paragraph 1: "The Tisdale RST sampling site ... 8-ft diameter cones."
# Collapses the italic "Oncorhynchus tshawytscha" / "O. mykiss" runs and the
# "two " run into the surrounding non-italic run's formatting.
Merge-Span "The Tisdale RST sampling site is operated by the California Department of Fish and Wildlife (CDFW) to obtain information on the temporal distribution, relative abundance, and composition of race and species of juvenile Chinook salmon (Oncorhynchus tshawytscha) and steelhead trout (O. mykiss) emigrating from the upper Sacramento River and tributaries to the Sacramento-San Joaquin Delta (Delta). The project collects data near the Tisdale Weir on the Sacramento River, using two paired rotary screw traps (RSTs) outfitted with two 8-ft diameter cones."

# Paragraph: "The RST monitoring site at Tisdale Weir was established meet a requirement ..."
Merge-Span "The RST monitoring site at Tisdale Weir was established meet a requirement of the 2011 amendment to the reasonable and prudent alternative (RPA) of the 2009 biological and conference opinion (BO) on the long-term operations of the Central Valley Project (CVP) and State Water Project (SWP). The amendment required the Bureau of Reclamation (USBR) and the California Department of Water Resources (DWR) to fund a new juvenile salmonid monitoring site on the Sacramento River between Red Bluff Diversion Dam (RBDD) and Knights Landing. The purpose of the new site was to provide early warning of fish movement and determine survival of listed fish species leaving spawning habitat in the upper Sacramento River.  "

# Paragraph: "CDFW issued ITP 2081-2019-066-00 ..." - merge "CDFW" + " issued ITP..." runs
Merge-Span "CDFW issued ITP 2081-2019-066-00 to DWR on March 31, 2020, for the long-term operation of the SWP in the Delta. "

# Same paragraph: "Condition 7.5.2 of the ITP requires ... Data from the Tisdale RST will be used ... JPE modeling approaches."
Merge-Span "Condition 7.5.2 of the ITP requires the development and establishment of a spring-run Chinook salmon juvenile production estimate (JPE) to increase understanding regarding the impacts water operations have on the spring-run Chinook salmon population in the Sacramento River watershed and inform the development of minimization measures to reduce take of spring-run Chinook salmon at Delta fish salvage facilities. Data from the Tisdale RST will be used along with other datasets from juvenile salmonid monitoring programs in the Sacramento River Watershed to inform the development of JPE modeling approaches."

# Paragraph: "... ) to understand the movement of juvenile salmon in the Sacramento River Watershed to estimate the number of winter-run and spring-run Chinook salmon that have entered the Delta. "
Merge-Span ") to understand the movement of juvenile salmon in the Sacramento River Watershed to estimate the number of winter-run and spring-run Chinook salmon that have entered the Delta. "

# Same paragraph: " is a real-time operations monitoring team required by Condition of Approval 8.1.2 of the ITP which meets weekly ... in the Delta."
Merge-Span " is a real-time operations monitoring team required by Condition of Approval 8.1.2 of the ITP which meets weekly from October through June, to provide advice for real-time management of SWP operations to DWR, CDFW, and the Water Operation Management Team (WOMT) to minimize take of winter-run and spring-run Chinook salmon in the Delta."

Write-Output "Done merging runs"
